$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 32 (A32) currently has the special "last row" date format (YYYY-MM-DD).
# It becomes a regular data row, so give it the same format as the other
# data rows (YYYY-MM-DD HH:MM:SS), matching A2:A31.
$ws.Range("A32").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Add the new last row (33) with the next day's values, using the special
# "last row" date format that A32 used to have.
$ws.Range("A33").Value = 45773
$ws.Range("A33").NumberFormat = "YYYY-MM-DD"
$ws.Range("B33").Value = 134
$ws.Range("C33").Value = 136
$ws.Range("D33").Value = 134
